$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 3517.647
$ws.Cells.Item(64, 9).Value = 3500
$ws.Cells.Item(64, 10).Value = 3600
$ws.Cells.Item(64, 11).Value = 3500
$ws.Cells.Item(64, 12).Value = 3600
$ws.Cells.Item(64, 13).Value = -3252
$ws.Cells.Item(64, 14).Value = -4096
$ws.Cells.Item(67, 8).Value = 3517.647
$ws.Cells.Item(67, 9).Value = 3500
$ws.Cells.Item(67, 10).Value = 3600
$ws.Cells.Item(67, 11).Value = 3500
$ws.Cells.Item(67, 12).Value = 3600
$ws.Cells.Item(67, 13).Value = -2642
$ws.Cells.Item(67, 14).Value = -5316
$ws.Cells.Item(106, 8).Value = 4600
$ws.Cells.Item(106, 9).Value = 2000
$ws.Cells.Item(106, 10).Value = 7200
$ws.Cells.Item(106, 11).Value = 2000
$ws.Cells.Item(106, 12).Value = 7200
$ws.Cells.Item(106, 13).Value = -1369
$ws.Cells.Item(106, 14).Value = -8462
$ws.Cells.Item(107, 8).Value = 1417.1904
$ws.Cells.Item(107, 9).Value = 1758.75
$ws.Cells.Item(107, 11).Value = 1758.75
$ws.Cells.Item(107, 13).Value = 161.25
$ws.Cells.Item(112, 8).Value = 25642244
$ws.Cells.Item(112, 9).Value = 250000460
$ws.Cells.Item(112, 10).Value = 1303.8
$ws.Cells.Item(112, 11).Value = 750001380
$ws.Cells.Item(112, 12).Value = 3911.4
$ws.Cells.Item(112, 13).Value = -750000272
$ws.Cells.Item(112, 14).Value = -6127.4
$ws.Cells.Item(137, 8).Value = 1135702.1
$ws.Cells.Item(137, 9).Value = 2071799
$ws.Cells.Item(137, 10).Value = 2532.3684
$ws.Cells.Item(137, 11).Value = 6215397
$ws.Cells.Item(137, 12).Value = 7597.1052
$ws.Cells.Item(137, 13).Value = -6212847
$ws.Cells.Item(137, 14).Value = -12697.1052
$ws.Cells.Item(138, 8).Value = 6340.96
$ws.Cells.Item(138, 9).Value = 997.35
$ws.Cells.Item(138, 10).Value = 7676.8623
$ws.Cells.Item(138, 11).Value = 2992.05
$ws.Cells.Item(138, 12).Value = 23030.5869
$ws.Cells.Item(138, 13).Value = 2147.95
$ws.Cells.Item(138, 14).Value = -33310.58689999999

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3466.976
$ws.Cells.Item(61, 9).Value = 1017.2258
$ws.Cells.Item(61, 10).Value = 10370.818
$ws.Cells.Item(61, 11).Value = 1017.2258
$ws.Cells.Item(61, 12).Value = 10370.818
$ws.Cells.Item(61, 13).Value = -805.2258
$ws.Cells.Item(61, 14).Value = -10794.818
$ws.Cells.Item(74, 8).Value = 5428.7915
$ws.Cells.Item(74, 9).Value = 8352.727999999999
$ws.Cells.Item(74, 10).Value = 2954.6924
$ws.Cells.Item(74, 11).Value = 8352.727999999999
$ws.Cells.Item(74, 12).Value = 2954.6924
$ws.Cells.Item(74, 13).Value = -7478.727999999999
$ws.Cells.Item(74, 14).Value = -4702.6924
$ws.Cells.Item(77, 8).Value = 5428.7915
$ws.Cells.Item(77, 9).Value = 8352.727999999999
$ws.Cells.Item(77, 10).Value = 2954.6924
$ws.Cells.Item(77, 11).Value = 41763.64
$ws.Cells.Item(77, 12).Value = 14773.462
$ws.Cells.Item(77, 13).Value = -37395.64
$ws.Cells.Item(77, 14).Value = -23509.462
$ws.Cells.Item(132, 8).Value = 1421.8219
$ws.Cells.Item(132, 9).Value = 882.617
$ws.Cells.Item(132, 11).Value = 2647.851
$ws.Cells.Item(132, 13).Value = -117.8509999999997
$ws.Cells.Item(136, 8).Value = 3466.976
$ws.Cells.Item(136, 9).Value = 1017.2258
$ws.Cells.Item(136, 10).Value = 10370.818
$ws.Cells.Item(136, 11).Value = 3051.6774
$ws.Cells.Item(136, 12).Value = 31112.454
$ws.Cells.Item(136, 13).Value = -501.6774
$ws.Cells.Item(136, 14).Value = -36212.454

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1386.1666
$ws.Cells.Item(107, 9).Value = 1275.4828
$ws.Cells.Item(107, 10).Value = 1844.7142
$ws.Cells.Item(107, 11).Value = 1275.4828
$ws.Cells.Item(107, 12).Value = 1844.7142
$ws.Cells.Item(107, 13).Value = 644.5172
$ws.Cells.Item(107, 14).Value = -5684.7142
$ws.Cells.Item(134, 8).Value = 4202.2446
$ws.Cells.Item(134, 9).Value = 1360.6471
$ws.Cells.Item(134, 11).Value = 4081.9413
$ws.Cells.Item(134, 13).Value = -1546.9413

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2745.8
$ws.Cells.Item(31, 9).Value = 1109.375
$ws.Cells.Item(31, 10).Value = 4616
$ws.Cells.Item(31, 11).Value = 1109.375
$ws.Cells.Item(31, 12).Value = 4616
$ws.Cells.Item(31, 13).Value = -814.375
$ws.Cells.Item(31, 14).Value = -5206
$ws.Cells.Item(34, 8).Value = 2745.8
$ws.Cells.Item(34, 9).Value = 1109.375
$ws.Cells.Item(34, 10).Value = 4616
$ws.Cells.Item(34, 11).Value = 1109.375
$ws.Cells.Item(34, 12).Value = 4616
$ws.Cells.Item(34, 13).Value = -907.375
$ws.Cells.Item(34, 14).Value = -5020
$ws.Cells.Item(58, 8).Value = 2507.9495
$ws.Cells.Item(58, 9).Value = 1562.5605
$ws.Cells.Item(58, 10).Value = 7307.615
$ws.Cells.Item(58, 11).Value = 1562.5605
$ws.Cells.Item(58, 12).Value = 7307.615
$ws.Cells.Item(58, 13).Value = -1359.5605
$ws.Cells.Item(58, 14).Value = -7713.615
$ws.Cells.Item(132, 8).Value = 3278.121
$ws.Cells.Item(132, 9).Value = 2792.1428
$ws.Cells.Item(132, 11).Value = 8376.428400000001
$ws.Cells.Item(132, 13).Value = -5846.428400000001
$ws.Cells.Item(134, 8).Value = 2013.125
$ws.Cells.Item(134, 9).Value = 1142.0588
$ws.Cells.Item(134, 10).Value = 4128.5713
$ws.Cells.Item(134, 11).Value = 3426.1764
$ws.Cells.Item(134, 12).Value = 12385.7139
$ws.Cells.Item(134, 13).Value = -891.1764000000003
$ws.Cells.Item(134, 14).Value = -17455.7139
$ws.Cells.Item(136, 8).Value = 2507.9495
$ws.Cells.Item(136, 9).Value = 1562.5605
$ws.Cells.Item(136, 10).Value = 7307.615
$ws.Cells.Item(136, 11).Value = 4687.681500000001
$ws.Cells.Item(136, 12).Value = 21922.845
$ws.Cells.Item(136, 13).Value = -2137.681500000001
$ws.Cells.Item(136, 14).Value = -27022.845

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 3254.257
$ws.Cells.Item(122, 9).Value = 1212
$ws.Cells.Item(122, 11).Value = 10908
$ws.Cells.Item(122, 13).Value = -8458

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3003.6843
$ws.Cells.Item(132, 9).Value = 2113.75
$ws.Cells.Item(132, 10).Value = 3241
$ws.Cells.Item(132, 11).Value = 6341.25
$ws.Cells.Item(132, 12).Value = 9723
$ws.Cells.Item(132, 13).Value = -3811.25
$ws.Cells.Item(132, 14).Value = -14783

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 2331.6667
$ws.Cells.Item(61, 9).Value = 2333.3333
$ws.Cells.Item(61, 11).Value = 2333.3333
$ws.Cells.Item(61, 13).Value = -2131.3333
$ws.Cells.Item(113, 8).Value = 2331.6667
$ws.Cells.Item(113, 9).Value = 2333.3333
$ws.Cells.Item(113, 11).Value = 2333.3333
$ws.Cells.Item(113, 13).Value = -163.3332999999998
$ws.Cells.Item(128, 8).Value = 41843.637
$ws.Cells.Item(128, 10).Value = 41843.637
$ws.Cells.Item(128, 12).Value = 41843.637
$ws.Cells.Item(128, 14).Value = -51803.637
$ws.Cells.Item(132, 8).Value = 28677.076
$ws.Cells.Item(132, 9).Value = 50000.8
$ws.Cells.Item(132, 11).Value = 150002.4
$ws.Cells.Item(132, 13).Value = -147472.4
$ws.Cells.Item(136, 8).Value = 4046.4443
$ws.Cells.Item(136, 9).Value = 1981.0769
$ws.Cells.Item(136, 10).Value = 5964.2856
$ws.Cells.Item(136, 11).Value = 5943.2307
$ws.Cells.Item(136, 12).Value = 17892.8568
$ws.Cells.Item(136, 13).Value = -3393.2307
$ws.Cells.Item(136, 14).Value = -22992.8568

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(131, 8).Value = 67903.336
$ws.Cells.Item(131, 10).Value = 67903.336
$ws.Cells.Item(131, 12).Value = 67903.336
$ws.Cells.Item(131, 14).Value = -77983.336
$ws.Cells.Item(132, 8).Value = 1745.2858
$ws.Cells.Item(132, 9).Value = 1127.7188
$ws.Cells.Item(132, 10).Value = 8332.666999999999
$ws.Cells.Item(132, 11).Value = 3383.1564
$ws.Cells.Item(132, 12).Value = 24998.001
$ws.Cells.Item(132, 13).Value = -853.1564000000003
$ws.Cells.Item(132, 14).Value = -30058.001
$ws.Cells.Item(136, 8).Value = 2963.919
$ws.Cells.Item(136, 9).Value = 1948.7142
$ws.Cells.Item(136, 10).Value = 6122.3335
$ws.Cells.Item(136, 11).Value = 5846.142599999999
$ws.Cells.Item(136, 12).Value = 18367.0005
$ws.Cells.Item(136, 13).Value = -3296.142599999999
$ws.Cells.Item(136, 14).Value = -23467.0005
